$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: Tagbar (a plugin to list functions)
$ws.Range("A31").Value = "Tagbar"
$ws.Range("B31").Value = "A plugin to list functions "
$ws.Range("C31").Value = "# Download and install with Pathogen`nNavigate to .vim/bundle and download: https://github.com/majutsushi/tagbar`n# Prequisition`nCtags, util that can be installed with dpkg/ apt-get/ apt-cyg"

# Match the row height used by the new entry
$ws.Rows.Item(31).RowHeight = 60

# Scroll the view down to the new row and move the selection past it,
# matching the author's final cursor position.
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 1
$ws.Range("A32").Select()
